$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.897.18"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "2.303.92"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.40"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.53"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -2.00%  "
$ws.Range("D9").Value = "2.302.71"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.100"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.54"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.331"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.31"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").Value = "59.889.51"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "2.715.56"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "2.311.39"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.46"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.07"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.35"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.54"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.61"
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("E25").Value = "  -2.57%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.71"
$ws.Range("E27").Value = "  -3.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.34"
$ws.Range("E28").Value = "  +3.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.60"
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.72"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("B31").Value = "SuiNetwork"
$ws.Range("C31").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").Value = "0.0$([char]0x2083)0723"
$ws.Range("E32").Value = "  -2.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.81"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("E34").Value = "  +2.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.377"
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.98"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "314.82"
$ws.Range("E40").Value = "  +3.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.11"
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.80"
$ws.Range("E43").Value = "  -3.59%  "
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("E46").Value = "  +2.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.73"
$ws.Range("E47").Value = "  +2.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0489"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").Value = "0.0$([char]0x2086)0219"
$ws.Range("E49").Value = "  +18.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0212"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("E51").Value = "  +0.16%  "
